# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker "JAVIER NARVAEZ ARAGON" (row 21) is removed from the statement,
# and the remaining "CARMELO YESITH RUIZ CONTRERAS" rows are re-sorted by
# Periodo Mora in ascending order (2108..2112 instead of 2112..2108).
# Totals (Valor Mora, Cant. Trabajadores, Cant. Periodos) are refreshed to
# reflect the updated data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Refresh summary figures for the updated data set ---
$ws.Range("E11").Value = 181705
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 5

# --- Re-sort the "Periodo Mora" column (E16:E20) ascending: 2108..2112 ---
$ws.Range("E16").Value = "2108"
$ws.Range("E17").Value = "2109"
$ws.Range("E18").Value = "2110"
$ws.Range("E19").Value = "2111"
$ws.Range("E20").Value = "2112"

# --- Drop the last worker's row (JAVIER NARVAEZ ARAGON / 1912), row 20 picks
#     up the "last row" border formatting that row 21 used to have ---
$ws.Range("B21:J21").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)
$ws.Rows(21).Delete()
